$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H98").Value = 9220
$ws.Range("I98").Value = 5821.6665
$ws.Range("K98").Value = 5821.6665
$ws.Range("M98").Value = -4323.6665
$ws.Range("H116").Value = 4127.8667
$ws.Range("I116").Value = 3699.6
$ws.Range("J116").Value = 4984.4
$ws.Range("K116").Value = 3699.6
$ws.Range("L116").Value = 4984.4
$ws.Range("M116").Value = -257.5999999999999
$ws.Range("N116").Value = -11868.4
$ws.Range("H122").Value = 9220
$ws.Range("I122").Value = 5821.6665
$ws.Range("K122").Value = 17464.9995
$ws.Range("M122").Value = -15014.9995
$ws.Range("H135").Value = 32258752
$ws.Range("I135").Value = 328.75
$ws.Range("K135").Value = 2958.75
$ws.Range("M135").Value = -423.75
$ws.Range("H138").Value = 1459.2858
$ws.Range("I138").Value = 644.8333
$ws.Range("J138").Value = 1932.1936
$ws.Range("K138").Value = 1934.4999
$ws.Range("L138").Value = 5796.5808
$ws.Range("M138").Value = 3205.5001
$ws.Range("N138").Value = -16076.5808

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3226.092
$ws.Range("I32").Value = 2883.1343
$ws.Range("J32").Value = 5779.222
$ws.Range("K32").Value = 2883.1343
$ws.Range("L32").Value = 5779.222
$ws.Range("M32").Value = -2596.1343
$ws.Range("N32").Value = -6353.222
$ws.Range("H45").Value = 1040.8788
$ws.Range("I45").Value = 985.95
$ws.Range("J45").Value = 1125.3846
$ws.Range("K45").Value = 985.95
$ws.Range("L45").Value = 1125.3846
$ws.Range("M45").Value = -608.95
$ws.Range("N45").Value = -1879.3846
$ws.Range("H61").Value = 1318.8055
$ws.Range("I61").Value = 1188.6897
$ws.Range("J61").Value = 1857.8572
$ws.Range("K61").Value = 1188.6897
$ws.Range("L61").Value = 1857.8572
$ws.Range("M61").Value = -976.6896999999999
$ws.Range("N61").Value = -2281.8572
$ws.Range("H74").Value = 962.9643
$ws.Range("I74").Value = 715.125
$ws.Range("K74").Value = 715.125
$ws.Range("M74").Value = 158.875
$ws.Range("H77").Value = 962.9643
$ws.Range("I77").Value = 715.125
$ws.Range("K77").Value = 3575.625
$ws.Range("M77").Value = 792.375
$ws.Range("H128").Value = 99990
$ws.Range("J128").Value = 99990
$ws.Range("L128").Value = 99990
$ws.Range("N128").Value = -109950
$ws.Range("H132").Value = 2263.7778
$ws.Range("I132").Value = 1910.7142
$ws.Range("K132").Value = 5732.142599999999
$ws.Range("M132").Value = -3202.142599999999
$ws.Range("H136").Value = 1318.8055
$ws.Range("I136").Value = 1188.6897
$ws.Range("J136").Value = 1857.8572
$ws.Range("K136").Value = 3566.0691
$ws.Range("L136").Value = 5573.571599999999
$ws.Range("M136").Value = -1016.0691
$ws.Range("N136").Value = -10673.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H99").Value = 45455824
$ws.Range("I99").Value = 50001210
$ws.Range("K99").Value = 50001210
$ws.Range("M99").Value = -49999712
$ws.Range("H134").Value = 4831.394
$ws.Range("I134").Value = 1225.6818
$ws.Range("J134").Value = 12042.818
$ws.Range("K134").Value = 3677.0454
$ws.Range("L134").Value = 36128.454
$ws.Range("M134").Value = -1142.0454
$ws.Range("N134").Value = -41198.454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1135.0469
$ws.Range("I31").Value = 1106.3276
$ws.Range("J31").Value = 1412.6666
$ws.Range("K31").Value = 1106.3276
$ws.Range("L31").Value = 1412.6666
$ws.Range("M31").Value = -811.3276000000001
$ws.Range("N31").Value = -2002.6666
$ws.Range("H34").Value = 1135.0469
$ws.Range("I34").Value = 1106.3276
$ws.Range("J34").Value = 1412.6666
$ws.Range("K34").Value = 1106.3276
$ws.Range("L34").Value = 1412.6666
$ws.Range("M34").Value = -904.3276000000001
$ws.Range("N34").Value = -1816.6666
$ws.Range("H94").Value = 1048.6666
$ws.Range("I94").Value = 789
$ws.Range("J94").Value = 1256.4
$ws.Range("K94").Value = 789
$ws.Range("L94").Value = 1256.4
$ws.Range("M94").Value = -338
$ws.Range("N94").Value = -2158.4
$ws.Range("H125").Value = 17571.428
$ws.Range("J125").Value = 17571.428
$ws.Range("L125").Value = 17571.428
$ws.Range("N125").Value = -22491.428
$ws.Range("H132").Value = 1968.8064
$ws.Range("I132").Value = 1561.4
$ws.Range("J132").Value = 3666.3333
$ws.Range("K132").Value = 4684.200000000001
$ws.Range("L132").Value = 10998.9999
$ws.Range("M132").Value = -2154.200000000001
$ws.Range("N132").Value = -16058.9999
$ws.Range("H134").Value = 977.5599999999999
$ws.Range("I134").Value = 782
$ws.Range("K134").Value = 2346
$ws.Range("M134").Value = 189

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 1500
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 1500
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 4500
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -4726
$ws.Range("H9").Value = 372.75
$ws.Range("J9").Value = 372.75
$ws.Range("L9").Value = 1118.25
$ws.Range("N9").Value = -1566.25
$ws.Range("H10").Value = 50
$ws.Range("I10").Value = 50
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 150
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -11
$ws.Range("N10").ClearContents()
$ws.Range("H11").Value = 1371.6666
$ws.Range("I11").Value = 2379.6667
$ws.Range("K11").Value = 7139.000100000001
$ws.Range("M11").Value = -6999.000100000001
$ws.Range("H12").Value = 67.60714
$ws.Range("I12").Value = 98.666664
$ws.Range("K12").Value = 295.999992
$ws.Range("M12").Value = -122.999992
$ws.Range("H13").Value = 316.5
$ws.Range("I13").Value = 179.8
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 539.4000000000001
$ws.Range("L13").Value = 3000
$ws.Range("M13").Value = -371.4000000000001
$ws.Range("N13").Value = -3336
$ws.Range("H15").Value = 650
$ws.Range("I15").Value = 766.6667
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 2300.0001
$ws.Range("L15").Value = 900
$ws.Range("M15").Value = -2160.0001
$ws.Range("N15").Value = -1180
$ws.Range("H17").Value = 580
$ws.Range("I17").Value = 216.66667
$ws.Range("K17").Value = 650.00001
$ws.Range("M17").Value = -481.00001
$ws.Range("H39").Value = 2076.074
$ws.Range("J39").Value = 1860.5834
$ws.Range("L39").Value = 5581.7502
$ws.Range("N39").Value = -6169.7502
$ws.Range("H44").Value = 2357.7144
$ws.Range("J44").Value = 2926
$ws.Range("L44").Value = 8778
$ws.Range("N44").Value = -9574
$ws.Range("H75").Value = 100
$ws.Range("I75").Value = 100
$ws.Range("K75").Value = 300
$ws.Range("M75").Value = 698
$ws.Range("H78").Value = 100
$ws.Range("I78").Value = 100
$ws.Range("K78").Value = 900
$ws.Range("M78").Value = 4092
$ws.Range("H118").Value = 500
$ws.Range("I118").Value = 500
$ws.Range("K118").Value = 1500
$ws.Range("M118").Value = -257
$ws.Range("H122").Value = 746.8333
$ws.Range("J122").Value = 726.4
$ws.Range("L122").Value = 6537.599999999999
$ws.Range("N122").Value = -11437.6
$ws.Range("H130").Value = 1858.8889
$ws.Range("J130").Value = 2032.8572
$ws.Range("L130").Value = 6098.571599999999
$ws.Range("N130").Value = -16138.5716
$ws.Range("H131").Value = 13160169
$ws.Range("I131").Value = 166667090
$ws.Range("J131").Value = 2433.0571
$ws.Range("K131").Value = 500001270
$ws.Range("L131").Value = 7299.1713
$ws.Range("M131").Value = -499996230
$ws.Range("N131").Value = -17379.1713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 12668.8
$ws.Range("J97").Value = 12668.8
$ws.Range("L97").Value = 12668.8
$ws.Range("N97").Value = -14650.8
$ws.Range("H114").Value = 51333
$ws.Range("J114").Value = 51333
$ws.Range("L114").Value = 51333
$ws.Range("N114").Value = -60011

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 11831.5
$ws.Range("J94").Value = 11831.5
$ws.Range("L94").Value = 11831.5
$ws.Range("N94").Value = -13633.5

